$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the empty inline-string placeholder at B49 (cell should become truly blank)
$ws.Range("B49").ClearContents()

# Row 50
$ws.Range("A50").Value = "05/01/2026 08:41:53"
$ws.Range("B50").Value = "05/01 08:34"
$ws.Range("C50").Value = "Metrópoles"
$ws.Range("D50").Value = "Secretário de Reformas Econômicas de Haddad deixa governo; veja motivo"
$ws.Range("E50").Value = "https://www.metropoles.com/brasil/secretario-de-reformas-economicas-de-haddad-deixa-governo-veja-motivo"
$ws.Range("F50").Value = "lula"
$ws.Range("G50").Value = "s Barbosa Pinto ocupava o cargo desde 2023 e comandou reformas importantes para o governo Lula"

# Row 51
$ws.Range("A51").Value = "05/01/2026 08:41:54"
$ws.Range("B51").Value = "05/01 08:31"
$ws.Range("C51").Value = "Folha de S.Paulo - Mercado - Principal"
$ws.Range("D51").Value = "Secretário de Reformas Econômicas da Fazenda deixa o governo"
$ws.Range("E51").Value = "https://redir.folha.com.br/redir/online/mercado/rss091/*https://www1.folha.uol.com.br/mercado/2026/01/secretario-de-reformas-economicas-da-fazenda-deixa-o-governo.shtml"
$ws.Range("F51").Value = "ministério da fazenda"
$ws.Range("G51").Value = "conômicas do &lt;a href=`"https://www1.folha.uol.com.br/folha-topicos/ministerio-da-fazenda/`"&gt;Ministério da Fazenda&lt;/a&gt;, Marcos Pinto, foi exonerado a pedido na última sexta-feira (2).`n&lt;a href=`"https://red"

# Row 52
$ws.Range("A52").Value = "05/01/2026 08:41:55"
$ws.Range("B52").Value = "05/01 08:30"
$ws.Range("C52").Value = "g1 > Política"
$ws.Range("D52").Value = "Em ano eleitoral, mercado financeiro projeta queda do juro, inflação no limite na meta e desaceleração do ritmo de alta PIB"
$ws.Range("E52").Value = "https://g1.globo.com/economia/noticia/2026/01/05/em-ano-eleitoral-mercado-financeiro-projeta-queda-do-juro-inflacao-no-limite-na-meta-e-desaceleracao-do-ritmo-de-alta-pib.ghtml"
$ws.Range("F52").Value = "inflação"
$ws.Range("G52").Value = "o, queda dos juros, desaceleração no ritmo de crescimento do Produto Interno Bruto (PIB), &lt;b&gt;inflação&lt;/b&gt; dentro dos limites do regime de metas e taxa de câmbio estável. `nAs projeções fazem parte"

# Row 53
$ws.Range("A53").Value = "05/01/2026 08:41:56"
$ws.Range("B53").Value = "05/01 08:14"
$ws.Range("C53").Value = "Metrópoles"
$ws.Range("D53").Value = "AL: acidente deixa seis mortos, incluindo criança e recém-nascido"
$ws.Range("E53").Value = "https://www.metropoles.com/brasil/al-acidente-deixa-seis-mortos-incluindo-crianca-e-recem-nascido"
$ws.Range("F53").Value = "câmara"
$ws.Range("G53").Value = "Acidente aconteceu na AL-220. Uma das vítimas é filho do presidente da Câmara de Vereadores de Piranhas (AL)"
